$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match the formatting of the existing header cell (H1) on the new
# header cells I1:J1 by copying its format (reuses the same style index
# instead of minting a new one).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# New header cells for columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data row values
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
